$wb = $excel.ActiveWorkbook

# --- Swap columns A (Team) and B (Driver) for rows 4-23 on the
# "2026 Drivers Points" and "2026 Drivers Price" sheets, so the layout
# becomes (Driver, Team) like the header row already says. ---
$sheetNames = @("2026 Drivers Points", "2026 Drivers Price")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    for ($r = 4; $r -le 23; $r++) {
        $cellA = $ws.Cells.Item($r, 1)
        $cellB = $ws.Cells.Item($r, 2)
        $a = $cellA.Value2
        $b = $cellB.Value2
        $cellA.Value = $b
        $cellB.Value = $a
    }
}

# --- Update the selection shown on "2026 Drivers Points" (not the
# active sheet) to A2:B23. ---
$wsPoints = $wb.Worksheets.Item("2026 Drivers Points")
$wsPoints.Range("A2:B23").Select() | Out-Null

# --- Make "2026 Drivers Price" the active sheet and select A2:B23
# there too (this also moves tabSelected / workbookView.activeTab). ---
$wsPrice = $wb.Worksheets.Item("2026 Drivers Price")
$wsPrice.Activate()
$wsPrice.Range("A2:B23").Select() | Out-Null
